# Adding new manifesto runner script
# - label the two data blocks in column A ("EDI" for the raw rates, "SCOT (Scaled)"
#   for the scaled rates) and resize column A to fit the new labels
# - leave the selection where the author left it when the workbook was saved

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "EDI"
$ws.Range("A12").Value = "SCOT (Scaled)"

$ws.Columns.Item(1).ColumnWidth = 11.83

[void]$ws.Range("K21").Select()
